# Update "想去人数" (number interested) figures in the "展览" and "全部类型" sheets.
$wb = $excel.ActiveWorkbook

$updates = @{
    3  = 7768
    6  = 44
    9  = 6105
    10 = 163
    11 = 16
    12 = 33
    13 = 1846
    16 = 845
    17 = 179
    18 = 5549
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
